$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation results (B,D,E,F,G,H,I,K,O columns) for rows 2-25
# for the 380 kV case
$ws.Range("B2").Value = 0.1230997642205836
$ws.Range("D2").Value = 0.02024965061229267
$ws.Range("E2").Value = 0.4246523451747066
$ws.Range("F2").Value = 0.3812909462845866
$ws.Range("G2").Value = 0.2327342981425176
$ws.Range("H2").Value = 0.4012823701756574
$ws.Range("I2").Value = 0.3595752720144993
$ws.Range("K2").Value = 0.88416834969604
$ws.Range("O2").Value = 1.182574226925311
$ws.Range("B3").Value = 0.1090946938934252
$ws.Range("D3").Value = 0.0177731365342666
$ws.Range("E3").Value = 0.3704713352441757
$ws.Range("F3").Value = 0.3790747268770076
$ws.Range("G3").Value = 0.2316442115419974
$ws.Range("H3").Value = 0.4047570462664893
$ws.Range("I3").Value = 0.3658923434216614
$ws.Range("K3").Value = 0.775112570019445
$ws.Range("O3").Value = 1.187268892890259
$ws.Range("B4").Value = 0.1004821413835799
$ws.Range("D4").Value = 0.01624512838375836
$ws.Range("E4").Value = 0.3372853010806125
$ws.Range("F4").Value = 0.3780659973587746
$ws.Range("G4").Value = 0.2312525546235094
$ws.Range("H4").Value = 0.4071490205321879
$ws.Range("I4").Value = 0.3700330356078361
$ws.Range("K4").Value = 0.7078420932071481
$ws.Range("O4").Value = 1.19127207758352
$ws.Range("B5").Value = 0.09696937143516493
$ws.Range("D5").Value = 0.01562063408052694
$ws.Range("E5").Value = 0.3237805993375389
$ws.Range("F5").Value = 0.3777432050851033
$ws.Range("G5").Value = 0.2311624152696723
$ws.Range("H5").Value = 0.4081886927399623
$ws.Range("I5").Value = 0.3717862655556115
$ws.Range("K5").Value = 0.6803526643384998
$ws.Range("O5").Value = 1.19318427154235
$ws.Range("B6").Value = 0.09638590062394314
$ws.Range("D6").Value = 0.01551682882309535
$ws.Range("E6").Value = 0.3215392481842372
$ws.Range("F6").Value = 0.3776949298473369
$ws.Range("G6").Value = 0.2311516327252576
$ws.Range("H6").Value = 0.4083652483613562
$ws.Range("I6").Value = 0.372081365702094
$ws.Range("K6").Value = 0.6757835168507142
$ws.Range("O6").Value = 1.193518726125575
$ws.Range("B7").Value = 0.1004347790566555
$ws.Range("D7").Value = 0.0162367135361805
$ws.Range("E7").Value = 0.3371030975839915
$ws.Range("F7").Value = 0.3780612869869202
$ws.Range("G7").Value = 0.2312510581882918
$ws.Range("H7").Value = 0.4071627791519745
$ws.Range("I7").Value = 0.3700564136159947
$ws.Range("K7").Value = 0.7074716672278498
$ws.Range("O7").Value = 1.191296730132791
$ws.Range("B8").Value = 0.118273756886893
$ws.Range("D8").Value = 0.01939731557582292
$ws.Range("E8").Value = 0.4059528901139942
$ws.Range("F8").Value = 0.3804535761673122
$ws.Range("G8").Value = 0.2323005925383654
$ws.Range("H8").Value = 0.4024267359521971
$ws.Range("I8").Value = 0.3616990151581891
$ws.Range("K8").Value = 0.8466312182547711
$ws.Range("O8").Value = 1.183959796000195
$ws.Range("B9").Value = 0.1531391710722119
$ws.Range("D9").Value = 0.02553461656236777
$ws.Range("E9").Value = 0.5416896561627595
$ws.Range("F9").Value = 0.3879502977075191
$ws.Range("G9").Value = 0.2365784196250118
$ws.Range("H9").Value = 0.3951943351990863
$ws.Range("I9").Value = 0.347389442638832
$ws.Range("K9").Value = 1.117001462872111
$ws.Range("O9").Value = 1.178507294025081
$ws.Range("B10").Value = 0.1786719350272961
$ws.Range("D10").Value = 0.03000478290307029
$ws.Range("E10").Value = 0.6419719615983155
$ws.Range("F10").Value = 0.3951866138249613
$ws.Range("G10").Value = 0.2410982246578612
$ws.Range("H10").Value = 0.391138963364952
$ws.Range("I10").Value = 0.3381443602223602
$ws.Range("K10").Value = 1.314039766974247
$ws.Range("O10").Value = 1.180010670235703
$ws.Range("B11").Value = 0.190267163993866
$ws.Range("D11").Value = 0.03202954505057676
$ws.Range("E11").Value = 0.6877394825874035
$ws.Range("F11").Value = 0.3988579138360535
$ws.Range("G11").Value = 0.2434585171801587
$ws.Range("H11").Value = 0.3895685057546387
$ws.Range("I11").Value = 0.3342140455841918
$ws.Range("K11").Value = 1.403316990293831
$ws.Range("O11").Value = 1.181904415458348
$ws.Range("B12").Value = 0.1946548893815532
$ws.Range("D12").Value = 0.03279497180017188
$ws.Range("E12").Value = 0.705093716841219
$ws.Range("F12").Value = 0.4003030012269591
$ws.Range("G12").Value = 0.244396443664499
$ws.Range("H12").Value = 0.3890133639814479
$ws.Range("I12").Value = 0.3327653509999244
$ws.Range("K12").Value = 1.437071223920839
$ws.Range("O12").Value = 1.182796563015984
$ws.Range("B13").Value = 0.1937100583815834
$ws.Range("D13").Value = 0.03263018236597759
$ws.Range("E13").Value = 0.7013551188518363
$ws.Range("F13").Value = 0.3999893320782206
$ws.Range("G13").Value = 0.2441924752962308
$ws.Range("H13").Value = 0.3891311627340741
$ws.Range("I13").Value = 0.333075590091763
$ws.Range("K13").Value = 1.429804037311442
$ws.Range("O13").Value = 1.182596621929463
$ws.Range("B14").Value = 0.190628209342151
$ws.Range("D14").Value = 0.03209254364512049
$ws.Range("E14").Value = 0.6891667537940407
$ws.Range("F14").Value = 0.3989757012565107
$ws.Range("G14").Value = 0.2435347941948578
$ws.Range("H14").Value = 0.3895220405338762
$ws.Range("I14").Value = 0.3340940661290208
$ws.Range("K14").Value = 1.406095046427993
$ws.Range("O14").Value = 1.181974299138687
$ws.Range("B15").Value = 0.1887400706548306
$ws.Range("D15").Value = 0.03176305274766378
$ws.Range("E15").Value = 0.681704082994699
$ws.Range("F15").Value = 0.398361973928921
$ws.Range("G15").Value = 0.2431377047662977
$ws.Range("H15").Value = 0.3897666190138125
$ws.Range("I15").Value = 0.3347230743323024
$ws.Range("K15").Value = 1.391565633331027
$ws.Range("O15").Value = 1.181615933418442
$ws.Range("B16").Value = 0.1779137367371106
$ws.Range("D16").Value = 0.02987227979672724
$ws.Range("E16").Value = 0.6389840815519392
$ws.Range("F16").Value = 0.3949543454033915
$ws.Range("G16").Value = 0.2409501288131253
$ws.Range("H16").Value = 0.3912471260837549
$ws.Range("I16").Value = 0.3384067590308657
$ws.Range("K16").Value = 1.308197958481003
$ws.Range("O16").Value = 1.179911344773103
$ws.Range("B17").Value = 0.1712668612308335
$ws.Range("D17").Value = 0.02871007785428503
$ws.Range("E17").Value = 0.6128160081837564
$ws.Range("F17").Value = 0.3929612643221816
$ws.Range("G17").Value = 0.239686314887507
$ws.Range("H17").Value = 0.3922257117930172
$ws.Range("I17").Value = 0.340737124174197
$ws.Range("K17").Value = 1.256961963586207
$ws.Range("O17").Value = 1.179176255433617
$ws.Range("B18").Value = 0.1674419125052964
$ws.Range("D18").Value = 0.02804079046050134
$ws.Range("E18").Value = 0.5977787087748823
$ws.Range("F18").Value = 0.3918505970294319
$ws.Range("G18").Value = 0.2389880129893243
$ws.Range("H18").Value = 0.3928143853470942
$ws.Range("I18").Value = 0.3421034005364039
$ws.Range("K18").Value = 1.227458929715567
$ws.Range("O18").Value = 1.178867265532801
$ws.Range("B19").Value = 0.1661465433096652
$ws.Range("D19").Value = 0.02781404201774507
$ws.Range("E19").Value = 0.5926896706506284
$ws.Range("F19").Value = 0.3914806671985787
$ws.Range("G19").Value = 0.2387564819562158
$ws.Range("H19").Value = 0.393018130989546
$ws.Range("I19").Value = 0.342570447450159
$ws.Range("K19").Value = 1.21746402182589
$ws.Range("O19").Value = 1.178782162490052
$ws.Range("B20").Value = 0.1719746252361176
$ws.Range("D20").Value = 0.02883388142192445
$ws.Range("E20").Value = 0.6156001950492822
$ws.Range("F20").Value = 0.3931697343393239
$ws.Range("G20").Value = 0.2398178861381126
$ws.Range("H20").Value = 0.3921188669560536
$ws.Range("I20").Value = 0.3404863707505736
$ws.Range("K20").Value = 1.262419597566634
$ws.Range("O20").Value = 1.179242719995671
$ws.Range("B21").Value = 0.1915335105714746
$ws.Range("D21").Value = 0.03225049703200966
$ws.Range("E21").Value = 0.6927461324072794
$ws.Range("F21").Value = 0.3992719384775256
$ws.Range("G21").Value = 0.2437267701617571
$ws.Range("H21").Value = 0.3894061558596604
$ws.Range("I21").Value = 0.3337938391542128
$ws.Range("K21").Value = 1.413060407922387
$ws.Range("O21").Value = 1.182152331848584
$ws.Range("B22").Value = 0.2042979569592376
$ws.Range("D22").Value = 0.03447581315928971
$ws.Range("E22").Value = 0.7433007152722411
$ws.Range("F22").Value = 0.4035798939727187
$ws.Range("G22").Value = 0.2465388847464567
$ws.Range("H22").Value = 0.3878638566172867
$ws.Range("I22").Value = 0.3296508974429972
$ws.Range("K22").Value = 1.511202452081193
$ws.Range("O22").Value = 1.185074613273827
$ws.Range("B23").Value = 0.1974871139285028
$ws.Range("D23").Value = 0.03328883556756068
$ws.Range("E23").Value = 0.7163058420107404
$ws.Range("F23").Value = 0.4012513020224802
$ws.Range("G23").Value = 0.2450143223953916
$ws.Range("H23").Value = 0.3886658747052536
$ws.Range("I23").Value = 0.3318409122295076
$ws.Range("K23").Value = 1.458851208579006
$ws.Range("O23").Value = 1.183421196482868
$ws.Range("B24").Value = 0.1716546562324908
$ws.Range("D24").Value = 0.02877791331975033
$ws.Range("E24").Value = 0.6143414422907938
$ws.Range("F24").Value = 0.3930753754756608
$ws.Range("G24").Value = 0.2397583147086664
$ws.Range("H24").Value = 0.3921670902999779
$ws.Range("I24").Value = 0.3405996537733316
$ws.Range("K24").Value = 1.259952346817045
$ws.Range("O24").Value = 1.179212317509382
$ws.Range("B25").Value = 0.1437208677840545
$ws.Range("D25").Value = 0.02388100728190778
$ws.Range("E25").Value = 0.5048799663215533
$ws.Range("F25").Value = 0.3856200105241641
$ws.Range("G25").Value = 0.2351810699485952
$ws.Range("H25").Value = 0.3969303251233995
$ws.Range("I25").Value = 0.3510379051535288
$ws.Range("K25").Value = 1.044135908283323
$ws.Range("O25").Value = 1.179019338748134
